$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 6.40767
$ws.Cells.Item(2, 8).Value = 19.22301
$ws.Cells.Item(2, 9).Value = 0.03429978795594129
$ws.Cells.Item(2, 10).Value = 0.03429978795594129
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 8.131233999999999
$ws.Cells.Item(2, 14).Value = 24.393702
$ws.Cells.Item(2, 15).Value = 0.02090995573015822
$ws.Cells.Item(2, 16).Value = 0.02090995573015823
$ws.Cells.Item(2, 17).Value = 52.10226416477999
$ws.Cells.Item(2, 18).Value = 468.9203774830199
$ws.Cells.Item(2, 19).Value = 0.0007172070477125466
$ws.Cells.Item(2, 20).Value = 0.0007172070477125467

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 6.40767
$ws.Cells.Item(3, 8).Value = 19.22301
$ws.Cells.Item(3, 9).Value = 0.03429978795594129
$ws.Cells.Item(3, 10).Value = 0.03429978795594129
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 243.3763986666667
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.625857000534647
$ws.Cells.Item(3, 16).Value = 0.6258570005346471
$ws.Cells.Item(3, 17).Value = 1559.47564844444
$ws.Cells.Item(3, 18).Value = 14035.28083599996
$ws.Cells.Item(3, 19).Value = 0.02146676240907983
$ws.Cells.Item(3, 20).Value = 0.02146676240907983

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 6.40767
$ws.Cells.Item(4, 8).Value = 19.22301
$ws.Cells.Item(4, 9).Value = 0.03429978795594129
$ws.Cells.Item(4, 10).Value = 0.03429978795594129
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 103.9426383333333
$ws.Cells.Item(4, 14).Value = 311.827915
$ws.Cells.Item(4, 15).Value = 0.2672947262403034
$ws.Cells.Item(4, 16).Value = 0.2672947262403035
$ws.Cells.Item(4, 17).Value = 666.03012536935
$ws.Cells.Item(4, 18).Value = 5994.27112832415
$ws.Cells.Item(4, 19).Value = 0.009168152431783783
$ws.Cells.Item(4, 20).Value = 0.009168152431783786

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 6.40767
$ws.Cells.Item(5, 8).Value = 19.22301
$ws.Cells.Item(5, 9).Value = 0.03429978795594129
$ws.Cells.Item(5, 10).Value = 0.03429978795594129
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 33.41874933333333
$ws.Cells.Item(5, 14).Value = 100.256248
$ws.Cells.Item(5, 15).Value = 0.08593831749489127
$ws.Cells.Item(5, 16).Value = 0.08593831749489128
$ws.Cells.Item(5, 17).Value = 214.13631754072
$ws.Cells.Item(5, 18).Value = 1927.22685786648
$ws.Cells.Item(5, 19).Value = 0.00294766606736513
$ws.Cells.Item(5, 20).Value = 0.002947666067365131

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 157.7959793333333
$ws.Cells.Item(6, 8).Value = 473.387938
$ws.Cells.Item(6, 9).Value = 0.8446703140819405
$ws.Cells.Item(6, 10).Value = 0.8446703140819404
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 8.131233999999999
$ws.Cells.Item(6, 14).Value = 24.393702
$ws.Cells.Item(6, 15).Value = 0.02090995573015822
$ws.Cells.Item(6, 16).Value = 0.02090995573015823
$ws.Cells.Item(6, 17).Value = 1283.076032218497
$ws.Cells.Item(6, 18).Value = 11547.68428996647
$ws.Cells.Item(6, 19).Value = 0.01766201887403222
$ws.Cells.Item(6, 20).Value = 0.01766201887403222

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 157.7959793333333
$ws.Cells.Item(7, 8).Value = 473.387938
$ws.Cells.Item(7, 9).Value = 0.8446703140819405
$ws.Cells.Item(7, 10).Value = 0.8446703140819404
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 243.3763986666667
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.625857000534647
$ws.Cells.Item(7, 16).Value = 0.6258570005346471
$ws.Cells.Item(7, 17).Value = 38403.81717422643
$ws.Cells.Item(7, 18).Value = 345634.3545680378
$ws.Cells.Item(7, 19).Value = 0.5286428292119815
$ws.Cells.Item(7, 20).Value = 0.5286428292119815

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 157.7959793333333
$ws.Cells.Item(8, 8).Value = 473.387938
$ws.Cells.Item(8, 9).Value = 0.8446703140819405
$ws.Cells.Item(8, 10).Value = 0.8446703140819404
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 103.9426383333333
$ws.Cells.Item(8, 14).Value = 311.827915
$ws.Cells.Item(8, 15).Value = 0.2672947262403034
$ws.Cells.Item(8, 16).Value = 0.2672947262403035
$ws.Cells.Item(8, 17).Value = 16401.73041029881
$ws.Cells.Item(8, 18).Value = 147615.5736926893
$ws.Cells.Item(8, 19).Value = 0.2257759203658434
$ws.Cells.Item(8, 20).Value = 0.2257759203658434

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 157.7959793333333
$ws.Cells.Item(9, 8).Value = 473.387938
$ws.Cells.Item(9, 9).Value = 0.8446703140819405
$ws.Cells.Item(9, 10).Value = 0.8446703140819404
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 33.41874933333333
$ws.Cells.Item(9, 14).Value = 100.256248
$ws.Cells.Item(9, 15).Value = 0.08593831749489127
$ws.Cells.Item(9, 16).Value = 0.08593831749489128
$ws.Cells.Item(9, 17).Value = 5273.344279148513
$ws.Cells.Item(9, 18).Value = 47460.09851233662
$ws.Cells.Item(9, 19).Value = 0.07258954563008332
$ws.Cells.Item(9, 20).Value = 0.07258954563008332

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.083191333333334
$ws.Cells.Item(10, 8).Value = 9.249574000000001
$ws.Cells.Item(10, 9).Value = 0.0165040972710719
$ws.Cells.Item(10, 10).Value = 0.01650409727107189
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.131233999999999
$ws.Cells.Item(10, 14).Value = 24.393702
$ws.Cells.Item(10, 15).Value = 0.02090995573015822
$ws.Cells.Item(10, 16).Value = 0.02090995573015823
$ws.Cells.Item(10, 17).Value = 25.07015019810533
$ws.Cells.Item(10, 18).Value = 225.631351782948
$ws.Cells.Item(10, 19).Value = 0.0003450999433043385
$ws.Cells.Item(10, 20).Value = 0.0003450999433043385

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3.083191333333334
$ws.Cells.Item(11, 8).Value = 9.249574000000001
$ws.Cells.Item(11, 9).Value = 0.0165040972710719
$ws.Cells.Item(11, 10).Value = 0.01650409727107189
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 243.3763986666667
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.625857000534647
$ws.Cells.Item(11, 16).Value = 0.6258570005346471
$ws.Cells.Item(11, 17).Value = 750.3760031069451
$ws.Cells.Item(11, 18).Value = 6753.384027962506
$ws.Cells.Item(11, 19).Value = 0.01032920481460511
$ws.Cells.Item(11, 20).Value = 0.01032920481460511

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 3.083191333333334
$ws.Cells.Item(12, 8).Value = 9.249574000000001
$ws.Cells.Item(12, 9).Value = 0.0165040972710719
$ws.Cells.Item(12, 10).Value = 0.01650409727107189
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 103.9426383333333
$ws.Cells.Item(12, 14).Value = 311.827915
$ws.Cells.Item(12, 15).Value = 0.2672947262403034
$ws.Cells.Item(12, 16).Value = 0.2672947262403035
$ws.Cells.Item(12, 17).Value = 320.4750416731345
$ws.Cells.Item(12, 18).Value = 2884.275375058211
$ws.Cells.Item(12, 19).Value = 0.004411458161914502
$ws.Cells.Item(12, 20).Value = 0.004411458161914501

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 3.083191333333334
$ws.Cells.Item(13, 8).Value = 9.249574000000001
$ws.Cells.Item(13, 9).Value = 0.0165040972710719
$ws.Cells.Item(13, 10).Value = 0.01650409727107189
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 33.41874933333333
$ws.Cells.Item(13, 14).Value = 100.256248
$ws.Cells.Item(13, 15).Value = 0.08593831749489127
$ws.Cells.Item(13, 16).Value = 0.08593831749489128
$ws.Cells.Item(13, 17).Value = 103.0363983153725
$ws.Cells.Item(13, 18).Value = 927.3275848383521
$ws.Cells.Item(13, 19).Value = 0.001418334351247945
$ws.Cells.Item(13, 20).Value = 0.001418334351247945

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 19.52685066666667
$ws.Cells.Item(14, 8).Value = 58.580552
$ws.Cells.Item(14, 9).Value = 0.1045258006910464
$ws.Cells.Item(14, 10).Value = 0.1045258006910464
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 8.131233999999999
$ws.Cells.Item(14, 14).Value = 24.393702
$ws.Cells.Item(14, 15).Value = 0.02090995573015822
$ws.Cells.Item(14, 16).Value = 0.02090995573015823
$ws.Cells.Item(14, 17).Value = 158.7773920537227
$ws.Cells.Item(14, 18).Value = 1428.996528483504
$ws.Cells.Item(14, 19).Value = 0.002185629865109123
$ws.Cells.Item(14, 20).Value = 0.002185629865109123

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 19.52685066666667
$ws.Cells.Item(15, 8).Value = 58.580552
$ws.Cells.Item(15, 9).Value = 0.1045258006910464
$ws.Cells.Item(15, 10).Value = 0.1045258006910464
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 243.3763986666667
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.625857000534647
$ws.Cells.Item(15, 16).Value = 0.6258570005346471
$ws.Cells.Item(15, 17).Value = 4752.374592555133
$ws.Cells.Item(15, 18).Value = 42771.3713329962
$ws.Cells.Item(15, 19).Value = 0.06541820409898066
$ws.Cells.Item(15, 20).Value = 0.06541820409898066

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 19.52685066666667
$ws.Cells.Item(16, 8).Value = 58.580552
$ws.Cells.Item(16, 9).Value = 0.1045258006910464
$ws.Cells.Item(16, 10).Value = 0.1045258006910464
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 103.9426383333333
$ws.Cells.Item(16, 14).Value = 311.827915
$ws.Cells.Item(16, 15).Value = 0.2672947262403034
$ws.Cells.Item(16, 16).Value = 0.2672947262403035
$ws.Cells.Item(16, 17).Value = 2029.672376634342
$ws.Cells.Item(16, 18).Value = 18267.05138970908
$ws.Cells.Item(16, 19).Value = 0.02793919528076178
$ws.Cells.Item(16, 20).Value = 0.02793919528076179

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 19.52685066666667
$ws.Cells.Item(17, 8).Value = 58.580552
$ws.Cells.Item(17, 9).Value = 0.1045258006910464
$ws.Cells.Item(17, 10).Value = 0.1045258006910464
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 33.41874933333333
$ws.Cells.Item(17, 14).Value = 100.256248
$ws.Cells.Item(17, 15).Value = 0.08593831749489127
$ws.Cells.Item(17, 16).Value = 0.08593831749489128
$ws.Cells.Item(17, 17).Value = 652.5629276987662
$ws.Cells.Item(17, 18).Value = 5873.066349288896
$ws.Cells.Item(17, 19).Value = 0.008982771446194875
$ws.Cells.Item(17, 20).Value = 0.008982771446194875
